$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 22.6500000000001
$ws.Range("G2").Value = 0.0003133741931310441
$ws.Range("H2").Value = 0.00461915000501239
$ws.Range("K2").Value = 4.436170125210158
$ws.Range("L2").Value = "[1.8102658591115093, 7.062074391308806]"
$ws.Range("M2").Value = 0.000981230800188193
$ws.Range("N2").Value = 0.001962461600376386
$ws.Range("O2").Value = -1.773631888513156
$ws.Range("P2").Value = "[-2.490632013656773, -1.05663176336954]"
$ws.Range("Q2").Value = 0.0000016879174882245
$ws.Range("R2").Value = 0.000003375834976448999
$ws.Range("S2").Value = 13.36438398535868
$ws.Range("T2").Value = "[11.843950566246992, 14.884817404470365]"
$ws.Range("W2").Value = 6.393693693693724
$ws.Range("X2").Value = 3.809009009009027
$ws.Range("Y2").Value = 8.97837837837842

# Row 3 updates
$ws.Range("E3").Value = 24.28000000000036
$ws.Range("G3").Value = 0.009977053265212477
$ws.Range("H3").Value = 0.0318189635536127
$ws.Range("K3").Value = 4.112497067198738
$ws.Range("L3").Value = "[1.1404460616295502, 7.084548072767927]"
$ws.Range("M3").Value = 0.00685578237760609
$ws.Range("N3").Value = 0.00685578237760609
$ws.Range("O3").Value = 1.616395018964118
$ws.Range("P3").Value = "[0.5723422051585008, 2.6604478327697354]"
$ws.Range("Q3").Value = 0.002525203919105268
$ws.Range("R3").Value = 0.002525203919105268
$ws.Range("S3").Value = 13.67370910234408
$ws.Range("T3").Value = "[11.804373860444377, 15.543044344243791]"
$ws.Range("W3").Value = 18.03379379379406
$ws.Range("X3").Value = 13.99927927927948
$ws.Range("Y3").Value = 22.06830830830863
